# Auto-generated edit script: updates FFXIV leve-profit market-price
# snapshot values (columns H-N) across several sheets, per the
# "update Sheets via scheduled runner" commit.
$wb = $excel.ActiveWorkbook

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2600.2104
$ws.Range("I86").Value = 2391.3076
$ws.Range("J86").Value = 3052.8333
$ws.Range("K86").Value = 2391.3076
$ws.Range("L86").Value = 3052.8333
$ws.Range("M86").Value = -1268.3076

# ALC row 88
$ws.Range("H88").Value = 11666.667
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 11666.667
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 11666.667
$ws.Range("N88").Value = -12478.667
$ws.Range("M88").ClearContents()

# ALC row 89
$ws.Range("H89").Value = 2600.2104
$ws.Range("I89").Value = 2391.3076
$ws.Range("J89").Value = 3052.8333
$ws.Range("K89").Value = 11956.538
$ws.Range("L89").Value = 15264.1665
$ws.Range("M89").Value = -6340.538

# ALC row 91
$ws.Range("H91").Value = 11666.667
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 11666.667
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 11666.667
$ws.Range("N91").Value = -14474.667
$ws.Range("M91").ClearContents()

# ALC row 98
$ws.Range("H98").Value = 520.6
$ws.Range("I98").Value = 473.625
$ws.Range("J98").Value = 708.5
$ws.Range("K98").Value = 473.625
$ws.Range("L98").Value = 708.5
$ws.Range("M98").Value = 1024.375
$ws.Range("N98").Value = -3704.5

# ALC row 122
$ws.Range("H122").Value = 520.6
$ws.Range("I122").Value = 473.625
$ws.Range("J122").Value = 708.5
$ws.Range("K122").Value = 1420.875
$ws.Range("L122").Value = 2125.5
$ws.Range("M122").Value = 1029.125
$ws.Range("N122").Value = -7025.5

# ALC row 129
$ws.Range("H129").Value = 3213.6047
$ws.Range("I129").Value = 11745.223
$ws.Range("J129").Value = 955.2353000000001
$ws.Range("K129").Value = 35235.669
$ws.Range("L129").Value = 2865.7059
$ws.Range("M129").Value = -30235.669
$ws.Range("N129").Value = -12865.7059

# ALC row 138
$ws.Range("H138").Value = 2818.7
$ws.Range("I138").Value = 1537.0714
$ws.Range("J138").Value = 3673.1191
$ws.Range("K138").Value = 4611.2142
$ws.Range("L138").Value = 11019.3573
$ws.Range("M138").Value = 528.7857999999997
$ws.Range("N138").Value = -21299.3573

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 34997.33
$ws.Range("I32").Value = 10524.055
$ws.Range("J32").Value = 259335.67
$ws.Range("K32").Value = 10524.055
$ws.Range("L32").Value = 259335.67
$ws.Range("M32").Value = -10237.055
$ws.Range("N32").Value = -259909.67

# ARM row 37
$ws.Range("H37").Value = 15000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 15000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 15000
$ws.Range("N37").Value = -15546

# ARM row 44
$ws.Range("H44").Value = 14980
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 14980
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 14980
$ws.Range("N44").Value = -15956

# ARM row 55
$ws.Range("H55").Value = 6341.25
$ws.Range("I55").Value = 9000
$ws.Range("J55").Value = 5961.4287
$ws.Range("K55").Value = 9000
$ws.Range("L55").Value = 5961.4287
$ws.Range("M55").Value = -8685
$ws.Range("N55").Value = -6591.4287

# ARM row 74
$ws.Range("H74").Value = 717.9048
$ws.Range("I74").Value = 614.8889
$ws.Range("J74").Value = 1336
$ws.Range("K74").Value = 614.8889
$ws.Range("L74").Value = 1336
$ws.Range("M74").Value = 259.1111
$ws.Range("N74").Value = -3084

# ARM row 77
$ws.Range("H77").Value = 717.9048
$ws.Range("I77").Value = 614.8889
$ws.Range("J77").Value = 1336
$ws.Range("K77").Value = 3074.4445
$ws.Range("L77").Value = 6680
$ws.Range("M77").Value = 1293.5555
$ws.Range("N77").Value = -15416

# ARM row 80
$ws.Range("H80").Value = 24019.23
$ws.Range("I80").Value = 20000
$ws.Range("J80").Value = 24354.166
$ws.Range("K80").Value = 20000
$ws.Range("L80").Value = 24354.166
$ws.Range("M80").Value = -19002
$ws.Range("N80").Value = -26350.166

# ARM row 83
$ws.Range("H83").Value = 24019.23
$ws.Range("I83").Value = 20000
$ws.Range("J83").Value = 24354.166
$ws.Range("K83").Value = 60000
$ws.Range("L83").Value = 73062.49800000001
$ws.Range("M83").Value = -55008
$ws.Range("N83").Value = -83046.49800000001

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 96728.63
$ws.Range("I20").Value = 105400.6
$ws.Range("J20").Value = 10009
$ws.Range("K20").Value = 105400.6
$ws.Range("L20").Value = 10009
$ws.Range("M20").Value = -105153.6
$ws.Range("N20").Value = -10503

# BSM row 130
$ws.Range("H130").Value = 43219.938
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 43219.938
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 43219.938
$ws.Range("N130").Value = -53259.938

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35578.7
$ws.Range("I31").Value = 787.8461
$ws.Range("J31").Value = 50654.734
$ws.Range("K31").Value = 787.8461
$ws.Range("L31").Value = 50654.734
$ws.Range("M31").Value = -492.8461
$ws.Range("N31").Value = -51244.734

# CRP row 34
$ws.Range("H34").Value = 35578.7
$ws.Range("I34").Value = 787.8461
$ws.Range("J34").Value = 50654.734
$ws.Range("K34").Value = 787.8461
$ws.Range("L34").Value = 50654.734
$ws.Range("M34").Value = -585.8461
$ws.Range("N34").Value = -51058.734

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 574
$ws.Range("I68").Value = 643.3333
$ws.Range("J68").Value = 470
$ws.Range("K68").Value = 1929.9999
$ws.Range("L68").Value = 1410
$ws.Range("M68").Value = -1118.9999
$ws.Range("N68").Value = -3032

# CUL row 71
$ws.Range("H71").Value = 574
$ws.Range("I71").Value = 643.3333
$ws.Range("J71").Value = 470
$ws.Range("K71").Value = 5789.9997
$ws.Range("L71").Value = 4230
$ws.Range("M71").Value = -1733.9997
$ws.Range("N71").Value = -12342

# CUL row 131
$ws.Range("H131").Value = 804234.5600000001
$ws.Range("I131").Value = 609.1667
$ws.Range("J131").Value = 993322.9
$ws.Range("K131").Value = 1827.5001
$ws.Range("L131").Value = 2979968.7
$ws.Range("M131").Value = 3212.4999
$ws.Range("N131").Value = -2990048.7

# CUL row 134
$ws.Range("H134").Value = 4272.032
$ws.Range("I134").Value = 2141.8667
$ws.Range("J134").Value = 6269.0625
$ws.Range("K134").Value = 6425.6001
$ws.Range("L134").Value = 18807.1875
$ws.Range("M134").Value = -1355.6001
$ws.Range("N134").Value = -28947.1875

# CUL row 137
$ws.Range("H137").Value = 4214925
$ws.Range("I137").Value = 73419.28999999999
$ws.Range("J137").Value = 10013032
$ws.Range("K137").Value = 220257.87
$ws.Range("L137").Value = 30039096
$ws.Range("M137").Value = -215157.87
$ws.Range("N137").Value = -30049296

# GSM row 69
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 58500
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 58500
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 58500
$ws.Range("N69").Value = -59998

# GSM row 70
$ws.Range("H70").Value = 61283.945
$ws.Range("I70").Value = 85167.84
$ws.Range("J70").Value = 7002.364
$ws.Range("K70").Value = 85167.84
$ws.Range("L70").Value = 7002.364
$ws.Range("M70").Value = -84897.84
$ws.Range("N70").Value = -7542.364

# GSM row 72
$ws.Range("H72").Value = 58500
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 58500
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 175500
$ws.Range("N72").Value = -182988

# GSM row 73
$ws.Range("H73").Value = 61283.945
$ws.Range("I73").Value = 85167.84
$ws.Range("J73").Value = 7002.364
$ws.Range("K73").Value = 85167.84
$ws.Range("L73").Value = 7002.364
$ws.Range("M73").Value = -84231.84
$ws.Range("N73").Value = -8874.364

# GSM row 132
$ws.Range("H132").Value = 5213
$ws.Range("I132").Value = 4224.875
$ws.Range("J132").Value = 6342.2856
$ws.Range("K132").Value = 12674.625
$ws.Range("L132").Value = 19026.8568
$ws.Range("M132").Value = -10144.625
$ws.Range("N132").Value = -24086.8568

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4118.1816
$ws.Range("I46").Value = 550
$ws.Range("J46").Value = 4911.1113
$ws.Range("K46").Value = 550
$ws.Range("L46").Value = 4911.1113
$ws.Range("M46").Value = -362
$ws.Range("N46").Value = -5287.1113

# LTW row 132
$ws.Range("H132").Value = 6010.1816
$ws.Range("I132").Value = 6456.4443
$ws.Range("J132").Value = 4002
$ws.Range("K132").Value = 19369.3329
$ws.Range("L132").Value = 12006
$ws.Range("M132").Value = -16839.3329
$ws.Range("N132").Value = -17066

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4070.7144
$ws.Range("I122").Value = 2098
$ws.Range("J122").Value = 9002.5
$ws.Range("K122").Value = 6294
$ws.Range("L122").Value = 27007.5
$ws.Range("M122").Value = -3844
$ws.Range("N122").Value = -31907.5

# WVR row 132
$ws.Range("H132").Value = 15210.363
$ws.Range("I132").Value = 12128
$ws.Range("J132").Value = 16971.715
$ws.Range("K132").Value = 36384
$ws.Range("L132").Value = 50915.145
$ws.Range("M132").Value = -33854
$ws.Range("N132").Value = -55975.145

# WVR row 136
$ws.Range("H136").Value = 14045.372
$ws.Range("I136").Value = 24345.072
$ws.Range("J136").Value = 4213.841
$ws.Range("K136").Value = 73035.216
$ws.Range("L136").Value = 12641.523
$ws.Range("M136").Value = -70485.216
$ws.Range("N136").Value = -17741.523
